$d = $word.ActiveDocument
$d.Content.Find.Execute("But Docker isn't the only folks developing container technology.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "But Docker isn't the only one developing container technology.", 2)
